# The workbook originally had a "birth_date" column (G) with dd/mm/yyyy
# formatted dates, followed by a "role" column (H). This edit removes the
# birth_date column entirely (as if the user selected column G's header and
# deleted it), which shifts the role column left into G.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$col = $ws.Columns.Item(7)
$col.Select()
$col.Delete()
